$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 505.9
$ws.Range("I33").Value = 488.7143
$ws.Range("K33").Value = 488.7143
$ws.Range("M33").Value = -259.7143

# Row 98
$ws.Range("H98").Value = 4019.6875
$ws.Range("I98").Value = 4869.923
$ws.Range("J98").Value = 335.33334
$ws.Range("K98").Value = 4869.923
$ws.Range("L98").Value = 335.33334
$ws.Range("M98").Value = -3371.923
$ws.Range("N98").Value = -3331.33334

# Row 116
$ws.Range("H116").Value = 2892.8333
$ws.Range("I116").Value = 1542
$ws.Range("J116").Value = 3857.7144
$ws.Range("K116").Value = 1542
$ws.Range("L116").Value = 3857.7144
$ws.Range("M116").Value = 1900
$ws.Range("N116").Value = -10741.7144

# Row 122
$ws.Range("H122").Value = 4019.6875
$ws.Range("I122").Value = 4869.923
$ws.Range("J122").Value = 335.33334
$ws.Range("K122").Value = 14609.769
$ws.Range("L122").Value = 1006.00002
$ws.Range("M122").Value = -12159.769
$ws.Range("N122").Value = -5906.00002

# Row 125
$ws.Range("H125").Value = 3825.75
$ws.Range("I125").Value = 970
$ws.Range("J125").Value = 4233.7144
$ws.Range("K125").Value = 8730
$ws.Range("L125").Value = 38103.4296
$ws.Range("M125").Value = -6270
$ws.Range("N125").Value = -43023.4296

# Row 129
$ws.Range("H129").Value = 803.76746
$ws.Range("J129").Value = 868.2432
$ws.Range("L129").Value = 2604.7296
$ws.Range("N129").Value = -12604.7296

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7558.075
$ws.Range("I32").Value = 5851.7593
$ws.Range("K32").Value = 5851.7593
$ws.Range("M32").Value = -5564.7593

# Row 45
$ws.Range("H45").Value = 1310.1111
$ws.Range("I45").Value = 1300.2
$ws.Range("J45").Value = 1322.5
$ws.Range("K45").Value = 1300.2
$ws.Range("L45").Value = 1322.5
$ws.Range("M45").Value = -923.2
$ws.Range("N45").Value = -2076.5

# Row 61
$ws.Range("H61").Value = 71430190
$ws.Range("I61").Value = 111112290
$ws.Range("J61").Value = 2412.8
$ws.Range("K61").Value = 111112290
$ws.Range("L61").Value = 2412.8
$ws.Range("M61").Value = -111112078
$ws.Range("N61").Value = -2836.8

# Row 97
$ws.Range("H97").Value = 709.5714
$ws.Range("I97").Value = 711.1667
$ws.Range("K97").Value = 711.1667
$ws.Range("M97").Value = -215.1667

# Row 135
$ws.Range("H135").Value = 22009.6
$ws.Range("J135").Value = 22009.6
$ws.Range("L135").Value = 22009.6
$ws.Range("N135").Value = -32149.6

# Row 136
$ws.Range("H136").Value = 71430190
$ws.Range("I136").Value = 111112290
$ws.Range("J136").Value = 2412.8
$ws.Range("K136").Value = 333336870
$ws.Range("L136").Value = 7238.400000000001
$ws.Range("M136").Value = -333334320
$ws.Range("N136").Value = -12338.4

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1618.4615
$ws.Range("I134").Value = 1392.625
$ws.Range("K134").Value = 4177.875
$ws.Range("M134").Value = -1642.875

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 200001700
$ws.Range("I16").Value = 200001700
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 200001700
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -200001413
$ws.Range("N16").ClearContents()

# Row 58
$ws.Range("H58").Value = 5035.5
$ws.Range("I58").Value = 1232.1
$ws.Range("J58").Value = 11374.5
$ws.Range("K58").Value = 1232.1
$ws.Range("L58").Value = 11374.5
$ws.Range("M58").Value = -1029.1
$ws.Range("N58").Value = -11780.5

# Row 99
$ws.Range("H99").Value = 1611.1052
$ws.Range("I99").Value = 1539.3572
$ws.Range("K99").Value = 1539.3572
$ws.Range("M99").Value = -41.35719999999992

# Row 105
$ws.Range("H105").Value = 813.8570999999999
$ws.Range("I105").Value = 782.8333
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 782.8333
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 964.1667
$ws.Range("N105").Value = -4494

# Row 107
$ws.Range("H107").Value = 873.94446
$ws.Range("I107").Value = 493.16666
$ws.Range("J107").Value = 1635.5
$ws.Range("K107").Value = 493.16666
$ws.Range("L107").Value = 1635.5
$ws.Range("M107").Value = 1426.83334
$ws.Range("N107").Value = -5475.5

# Row 109
$ws.Range("H109").Value = 9450.5
$ws.Range("J109").Value = 9450.5
$ws.Range("L109").Value = 9450.5
$ws.Range("N109").Value = -11530.5

# Row 113
$ws.Range("H113").Value = 200001700
$ws.Range("I113").Value = 200001700
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 200001700
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -199999530
$ws.Range("N113").ClearContents()

# Row 126
$ws.Range("H126").Value = 1611.1052
$ws.Range("I126").Value = 1539.3572
$ws.Range("K126").Value = 4618.071599999999
$ws.Range("M126").Value = -2148.071599999999

# Row 136
$ws.Range("H136").Value = 5035.5
$ws.Range("I136").Value = 1232.1
$ws.Range("J136").Value = 11374.5
$ws.Range("K136").Value = 3696.3
$ws.Range("L136").Value = 34123.5
$ws.Range("M136").Value = -1146.3
$ws.Range("N136").Value = -39223.5

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1060
$ws.Range("I68").Value = 1450
$ws.Range("K68").Value = 4350
$ws.Range("M68").Value = -3539

# Row 71
$ws.Range("H71").Value = 1060
$ws.Range("I71").Value = 1450
$ws.Range("K71").Value = 13050
$ws.Range("M71").Value = -8994

# Row 82
$ws.Range("H82").Value = 8400.857
$ws.Range("I82").Value = 2522.6
$ws.Range("K82").Value = 7567.799999999999
$ws.Range("M82").Value = -7161.799999999999

# Row 85
$ws.Range("H85").Value = 8400.857
$ws.Range("I85").Value = 2522.6
$ws.Range("K85").Value = 7567.799999999999
$ws.Range("M85").Value = -6163.799999999999

# Row 98
$ws.Range("H98").Value = 1062.6923
$ws.Range("I98").Value = 1663.7142
$ws.Range("J98").Value = 361.5
$ws.Range("K98").Value = 4991.142599999999
$ws.Range("L98").Value = 1084.5
$ws.Range("M98").Value = -3493.142599999999
$ws.Range("N98").Value = -4080.5

# Row 122
$ws.Range("H122").Value = 1472.2354
$ws.Range("J122").Value = 2189.111
$ws.Range("L122").Value = 19701.999
$ws.Range("N122").Value = -24601.999

# Row 131
$ws.Range("H131").Value = 21309284
$ws.Range("I131").Value = 71428984
$ws.Range("J131").Value = 46382.848
$ws.Range("K131").Value = 214286952
$ws.Range("L131").Value = 139148.544
$ws.Range("M131").Value = -214281912
$ws.Range("N131").Value = -149228.544

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1537.1111
$ws.Range("I102").Value = 1610
$ws.Range("K102").Value = 1610
$ws.Range("M102").Value = 12

# Row 122
$ws.Range("H122").Value = 3115.3333
$ws.Range("I122").Value = 3418.5334
$ws.Range("J122").Value = 1599.3334
$ws.Range("K122").Value = 10255.6002
$ws.Range("L122").Value = 4798.0002
$ws.Range("M122").Value = -7805.600199999999
$ws.Range("N122").Value = -9698.0002

# Row 132
$ws.Range("H132").Value = 2823.5144
$ws.Range("I132").Value = 2647.3044
$ws.Range("J132").Value = 3161.25
$ws.Range("K132").Value = 7941.9132
$ws.Range("L132").Value = 9483.75
$ws.Range("M132").Value = -5411.9132
$ws.Range("N132").Value = -14543.75

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4809.615
$ws.Range("I40").Value = 2334.0908
$ws.Range("J40").Value = 18425
$ws.Range("K40").Value = 2334.0908
$ws.Range("L40").Value = 18425
$ws.Range("M40").Value = -2198.0908
$ws.Range("N40").Value = -18697

# Row 61
$ws.Range("H61").Value = 1257.6923
$ws.Range("J61").Value = 1526.75
$ws.Range("L61").Value = 1526.75
$ws.Range("N61").Value = -1930.75

# Row 100
$ws.Range("H100").Value = 1220
$ws.Range("I100").Value = 1033.3334
$ws.Range("K100").Value = 1033.3334
$ws.Range("M100").Value = -492.3334

# Row 113
$ws.Range("H113").Value = 1257.6923
$ws.Range("J113").Value = 1526.75
$ws.Range("L113").Value = 1526.75
$ws.Range("N113").Value = -5866.75

# Row 122
$ws.Range("H122").Value = 13159629
$ws.Range("I122").Value = 16668330
$ws.Range("K122").Value = 50004990
$ws.Range("M122").Value = -50002540

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 466.33334
$ws.Range("I107").Value = 392
$ws.Range("J107").Value = 512.0769
$ws.Range("K107").Value = 1176
$ws.Range("L107").Value = 1536.2307
$ws.Range("M107").Value = 744
$ws.Range("N107").Value = -5376.2307

# Row 122
$ws.Range("H122").Value = 12501302
$ws.Range("I122").Value = 13890209
$ws.Range("K122").Value = 41670627
$ws.Range("M122").Value = -41668177

# Row 126
$ws.Range("H126").Value = 71429790
$ws.Range("I126").Value = 111111840
$ws.Range("J126").Value = 2102
$ws.Range("K126").Value = 333335520
$ws.Range("L126").Value = 6306
$ws.Range("M126").Value = -333333050
$ws.Range("N126").Value = -11246
